$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block -----------------------------------------------------
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must stay TEXT (it already round-trips
# beyond Excel's 15-significant-digit numeric precision, and the source file
# stores it as an inline string). Writing it straight into Value would make
# Excel "smart type" it into a Double and round the last digit. So stage the
# text in a scratch cell that is explicitly formatted as Text, then copy/
# paste-special just the value into B3 - this keeps B3's own style (s="8")
# completely untouched, exactly like the target.
$scratch = $ws.Range("ZZ1000")
$scratch.NumberFormat = "@"
$scratch.Value = "2570314725427075"
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ----------------------------------------------
$ws.Range("D5").Value = "KONTOSTAND AM 18.06.2025"

# --- Transaction rows ----------------------------------------------------
$ws.Range("B6").Value = "21.06."
$ws.Range("C6").Value = "22.06."
$ws.Range("D6").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E6").Value = "57,30-"

$ws.Range("B7").Value = "23.06."
$ws.Range("C7").Value = "24.06."
$ws.Range("D7").Value = "BURGER KING Eckernförde"
$ws.Range("E7").Value = "29,22-"

$ws.Range("B8").Value = "27.06."
$ws.Range("C8").Value = "28.06."
$ws.Range("D8").Value = "ZALANDO MKTPLC EU YCEJDO"
$ws.Range("E8").Value = "43,30-"

$ws.Range("B9").Value = "28.06."
$ws.Range("C9").Value = "29.06."
$ws.Range("D9").Value = "KARTENZ./28.06 LIDL RO"
$ws.Range("E9").Value = "117,91-"

$ws.Range("B10").Value = "01.07."
$ws.Range("C10").Value = "02.07."
$ws.Range("D10").Value = "MCDONALDS Bersenbrück"
$ws.Range("E10").Value = "32,94-"

$ws.Range("B11").Value = "03.07."
$ws.Range("C11").Value = "04.07."
$ws.Range("D11").Value = "BEITRAG Allianz SE K-72013808"
$ws.Range("E11").Value = "56,86-"

# --- Closing balance line ------------------------------------------------
$ws.Range("D12").Value = "KONTOSTAND AM 07.07.2025"
$ws.Range("E12").Value = "337,53-"

# --- Next statement date --------------------------------------------------
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 14.07.2025"
